{"js": "// Version-control table: update the reviewer / approver names for the two\n// most recent revision rows (v2.8.1 and v2.4.1).\n//\n//   v2.8.1 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA)\" -> \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM)\"\n//               \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08:      \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\"  -> \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"\n//   v2.4.1 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA)\" -> \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"\n//               \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08:      \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\"  -> \"\u0e27\u0e23\u0e34\u0e28\u0e23\u0e32 (D)\"\n//\n// Each name/role pair is stored as two separate runs inside the same table\n// cell (one run for the name, one for the \" (ROLE)\" suffix), so each part\n// is replaced independently (searched for and replaced in place) in order\n// to keep the existing run-level formatting intact.\n\nconst table = context.document.body.tables.getFirst();\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Replace the first run of text `oldText` found inside `cell`'s body with\n// `newText`, keeping that run's own formatting (rFonts/sz/etc.) untouched.\nasync function replaceInCell(cell, oldText, newText) {\n  const found = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(`Text not found in cell: ${oldText}`);\n  }\n  found.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Row index 5 == version 2.8.1 (row 0 = header \"\u0e0a\u0e37\u0e48\u0e2d\u0e40\u0e2d\u0e01\u0e2a\u0e32\u0e23\", row 4 = column\n// headers \"\u0e40\u0e27\u0e2d\u0e23\u0e4c\u0e0a\u0e31\u0e19/\u0e27\u0e31\u0e19\u0e17\u0e35\u0e48/\u0e23\u0e32\u0e22\u0e25\u0e30\u0e40\u0e2d\u0e35\u0e22\u0e14/\u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a/\u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08\").\nconst row281 = table.rows.items[5];\nrow281.cells.load(\"items\");\nawait context.sync();\n\nconst responsible281 = row281.cells.items[3]; // \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a\nconst reviewer281 = row281.cells.items[4];    // \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08\n\nawait replaceInCell(responsible281, \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\", \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\");\nawait replaceInCell(responsible281, \" (QA)\", \" (DM)\");\n\nawait replaceInCell(reviewer281, \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\");\nawait replaceInCell(reviewer281, \"(SP)\", \" (TL)\");\n\n// Row index 6 == version 2.4.1.\nconst row241 = table.rows.items[6];\nrow241.cells.load(\"items\");\nawait context.sync();\n\nconst responsible241 = row241.cells.items[3]; // \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a\nconst reviewer241 = row241.cells.items[4];    // \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08\n\nawait replaceInCell(responsible241, \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\", \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\");\nawait replaceInCell(responsible241, \" (QA)\", \" (TL)\");\n\nawait replaceInCell(reviewer241, \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", \"\u0e27\u0e23\u0e34\u0e28\u0e23\u0e32\");\nawait replaceInCell(reviewer241, \"(SP)\", \" (D)\");\n", "ps1": "# Version-control table: update the reviewer / approver names for the two\n# most recent revision rows (v2.8.1 and v2.4.1).\n#\n#   v2.8.1 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA)\" -> \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM)\"\n#               \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08:      \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\"  -> \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"\n#   v2.4.1 row -> \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA)\" -> \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"\n#               \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08:      \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\"  -> \"\u0e27\u0e23\u0e34\u0e28\u0e23\u0e32 (D)\"\n#\n# Each name/role pair lives in two separate runs inside the same table cell\n# (one run for the name, one for the \" (ROLE)\" suffix). Restricting\n# Find.Execute to each cell's own Range keeps the replace scoped to that\n# single cell (instead of touching every other row that repeats the same\n# original text) and preserves each run's existing formatting.\n\nfunction Replace-InRange($range, [string]$findText, [string]$replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap:=0 (wdFindStop), Replace:=1 (wdReplaceOne) keeps the match (and the\n    # replacement) confined to the supplied range.\n    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1)\n    if (-not $ok) {\n        throw \"Find.Execute could not find '$findText' in the given range\"\n    }\n}\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Row 6 (1-based) == version 2.8.1 (row 1 = \"\u0e0a\u0e37\u0e48\u0e2d\u0e40\u0e2d\u0e01\u0e2a\u0e32\u0e23\", row 5 = column\n# headers \"\u0e40\u0e27\u0e2d\u0e23\u0e4c\u0e0a\u0e31\u0e19/\u0e27\u0e31\u0e19\u0e17\u0e35\u0e48/\u0e23\u0e32\u0e22\u0e25\u0e30\u0e40\u0e2d\u0e35\u0e22\u0e14/\u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a/\u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08\").\n$responsible281 = $table.Cell(6, 4)   # \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a\n$reviewer281 = $table.Cell(6, 5)      # \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08\n\nReplace-InRange $responsible281.Range \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\" \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\"\nReplace-InRange $responsible281.Range \" (QA)\" \" (DM)\"\n\nReplace-InRange $reviewer281.Range \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\nReplace-InRange $reviewer281.Range \"(SP)\" \" (TL)\"\n\n# Row 7 (1-based) == version 2.4.1.\n$responsible241 = $table.Cell(7, 4)   # \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a\n$reviewer241 = $table.Cell(7, 5)      # \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08\n\nReplace-InRange $responsible241.Range \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\nReplace-InRange $responsible241.Range \" (QA)\" \" (TL)\"\n\nReplace-InRange $reviewer241.Range \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" \"\u0e27\u0e23\u0e34\u0e28\u0e23\u0e32\"\nReplace-InRange $reviewer241.Range \"(SP)\" \" (D)\"\n"}
